$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)
$hdr = $sec.Headers.Item(1)
$f = $hdr.Range.Fields.Item(1)
$f.Delete()
Write-Output "done"
